$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must be forced to Text format
# so Excel stores them as text (matching the source data), not as numeric values.
$textCells = @(
    "D4",
    "D5",
    "D7",
    "D8",
    "D9",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D18",
    "D19",
    "D20",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '26.598.85'
$ws.Range('E2').Value = '  +4.08%  '
$ws.Range('D3').Value = '1.744.43'
$ws.Range('E3').Value = '  +4.53%  '
$ws.Range('D4').Value = '0.9992'
$ws.Range('D5').Value = '246.34'
$ws.Range('E5').Value = '  +3.86%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = '0.4823'
$ws.Range('E7').Value = '  +1.28%  '
$ws.Range('D8').Value = '0.2698'
$ws.Range('E8').Value = '  +3.53%  '
$ws.Range('D9').Value = '0.06260'
$ws.Range('E9').Value = '  +1.39%  '
$ws.Range('D10').Value = '1.744.22'
$ws.Range('E10').Value = '  +4.55%  '
$ws.Range('D11').Value = '0.07136'
$ws.Range('E11').Value = '  +1.80%  '
$ws.Range('D12').Value = '15.92'
$ws.Range('E12').Value = '  +7.72%  '
$ws.Range('D13').Value = '0.6255'
$ws.Range('E13').Value = '  +6.90%  '
$ws.Range('D14').Value = '4.527'
$ws.Range('E14').Value = '  +3.75%  '
$ws.Range('D15').Value = '77.51'
$ws.Range('D16').Value = '0.9999'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '26.598.17'
$ws.Range('D18').Value = '0.9998'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('D19').Value = '0.000006921'
$ws.Range('E19').Value = '  +2.92%  '
$ws.Range('D20').Value = '11.74'
$ws.Range('E20').Value = '  +2.86%  '
$ws.Range('D21').Value = '1.968.07'
$ws.Range('E21').Value = '  +4.50%  '
$ws.Range('D22').Value = '4.628'
$ws.Range('E22').Value = '  +4.39%  '
$ws.Range('D23').Value = '8.882'
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').Value = '5.374'
$ws.Range('E24').Value = '  +2.76%  '
$ws.Range('D25').Value = '136.34'
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('D26').Value = '15.38'
$ws.Range('E26').Value = '  +2.57%  '
$ws.Range('E27').Value = '  +6.24%  '
$ws.Range('D28').Value = '1.432'
$ws.Range('E28').Value = '  +3.72%  '
$ws.Range('D29').Value = '106.81'
$ws.Range('E29').Value = '  +2.26%  '
$ws.Range('D30').Value = '4.021'
$ws.Range('E30').Value = '  +0.78%  '
$ws.Range('D31').Value = '3.743'
$ws.Range('E31').Value = '  +3.47%  '
$ws.Range('E32').Value = '  +0.70%  '
$ws.Range('D33').Value = '0.04599'
$ws.Range('E33').Value = '  +6.36%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '2.616'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '0.6421'
$ws.Range('E35').Value = '  +6.02%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '0.9994'
$ws.Range('E36').Value = '  +4.84%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '0.9332'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('B38').Value = 'Quant'
$ws.Range('C38').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D38').Value = '114.00'
$ws.Range('E38').Value = '  +14.94%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.444'
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('D40').Value = '1.992'
$ws.Range('E40').Value = '  +7.65%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '1.004'
$ws.Range('E41').Value = '  +0.40%  '
$ws.Range('D42').Value = '5.794'
$ws.Range('E42').Value = '  +18.48%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.01512'
$ws.Range('E43').Value = '  +2.43%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.3923'
$ws.Range('E44').Value = '  +4.75%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = '0.1221'
$ws.Range('E45').Value = '  +9.58%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '6.767'
$ws.Range('E46').Value = '  +9.30%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.05338'
$ws.Range('E47').Value = '  +1.43%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '7.959'
$ws.Range('E48').Value = '  +6.31%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '30.83'
$ws.Range('E49').Value = '  +3.21%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '1.263'
$ws.Range('E50').Value = '  +4.98%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').Value = '0.3454'
$ws.Range('E51').Value = '  +3.84%  '
